$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 0) Remove the old "_GoBack" bookmark that precedes the manual page break,
#    before a new one is (re-)created elsewhere below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1) "...Gestión de Alumnosv1.0" paragraph: append " (preliminar)" + "."
#    runs, and move the "_GoBack" bookmark here.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found = $rng1.Find.Execute("Gestión de Alumnosv1.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Gestión de Alumnosv1.0' paragraph" }
$p1 = $rng1.Paragraphs(1).Range
$xml1 = '<w:p ' + $wNs + ' w:rsidR="00ED5CEE" w:rsidRDefault="00ED5CEE" w:rsidP="00ED5CEE">' + `
  '<w:pPr><w:keepLines/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0"/><w:ind w:left="734" w:hanging="425"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>C</w:t></w:r>' + `
  '<w:r w:rsidRPr="00ED5CEE"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">asos de uso de </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Gestión de Alumnosv1.0</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> (preliminar)</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Casos de uso de Gestión de Docentesv1.0" paragraph: append
#    " (preliminar)" + "." runs.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found = $rng2.Find.Execute("Casos de uso de Gestión de Docentesv1.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Casos de uso de Gestión de Docentesv1.0' paragraph" }
$p2 = $rng2.Paragraphs(1).Range
$xml2 = '<w:p ' + $wNs + ' w:rsidR="00ED5CEE" w:rsidRDefault="00ED5CEE" w:rsidP="00ED5CEE">' + `
  '<w:pPr><w:keepLines/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0"/><w:ind w:left="734" w:hanging="425"/><w:contextualSpacing/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Casos de uso de Gestión de Docentesv1.0</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> (preliminar)</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r>' + `
  '</w:p>'
$p2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) "Código fuente de gestión de perfiles." paragraph: "." -> "(preliminar)."
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found = $rng3.Find.Execute("Código fuente de gestión de perfiles.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Código fuente de gestión de perfiles.' paragraph" }
$p3 = $rng3.Paragraphs(1).Range
$xml3 = '<w:p ' + $wNs + ' w:rsidR="00ED5CEE" w:rsidRDefault="00ED5CEE" w:rsidP="00ED5CEE">' + `
  '<w:pPr><w:keepLines/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0"/><w:ind w:left="734" w:hanging="425"/><w:jc w:val="left"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="00596A49"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Código fuente de </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>gestión de perfiles</w:t></w:r>' + `
  '<w:r w:rsidRPr="00596A49"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>(preliminar).</w:t></w:r>' + `
  '</w:p>'
$p3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) "Documento Cierre del Proyecto" run: add <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$found = $rng4.Find.Execute("Documento Cierre del Proyecto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Documento Cierre del Proyecto' paragraph" }
$p4 = $rng4.Paragraphs(1).Range
$xml4 = '<w:p ' + $wNs + ' w:rsidR="00E66969" w:rsidRPr="00E726A0" w:rsidRDefault="00A1205C" w:rsidP="00E726A0">' + `
  '<w:pPr><w:keepLines/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0"/><w:ind w:left="714" w:hanging="357"/><w:jc w:val="left"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="00596A49"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Documento Cierre del Proyecto</w:t></w:r>' + `
  '</w:p>'
$p4.InsertXML($xml4)

# ---------------------------------------------------------------------------
# 5) Footer page-number field: cached result "3" -> "1".
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null

Write-Output "done"
